{"js": "// Insert a new numbered list item \"How to start with style\" right after\n// the existing \"Make counter\" list item, matching the same list\n// formatting (ListParagraph style, ilvl=0, numId=1).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Make counter\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error('Could not find paragraph with text \"Make counter\"');\n}\n\n// insertParagraph(\"After\") creates the new paragraph immediately after\n// the target, inheriting its paragraph properties (list style/numbering),\n// then we set the new paragraph's text.\ntarget.insertParagraph(\"How to start with style\", \"After\");\nawait context.sync();\n", "ps1": "# Insert a new numbered list item \"How to start with style\" right after\n# the existing \"Make counter\" list item, matching the same list\n# formatting (ListParagraph style, ilvl=0, numId=1).\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"Make counter\")\nif (-not $found) {\n    throw 'Could not find paragraph with text \"Make counter\"'\n}\n\n# $rng now spans the matched text; Paragraphs(1) is its containing paragraph.\n$target = $rng.Paragraphs(1)\n\n# InsertParagraphAfter() splits in a new paragraph right after $target,\n# inheriting $target's paragraph formatting (style/list numbering).\n$target.Range.InsertParagraphAfter()\n\n# The freshly inserted paragraph is now the target's successor; give it text.\n$newPara = $target.Next()\n$newPara.Range.Text = \"How to start with style\"\n"}
